$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Omar) ---
$ws.Range("D2").Value = 1276119111

# --- Update row 3 (was Reem Adel, now Gaber Moahmed / Gaber Mohamed / Gaber123) ---
$ws.Range("B3").Value = "Gaber Moahmed"
$ws.Range("E3").Value = "Gaber123"

# --- Clear rows 4-8 (remove the other students entirely, leaving a blank styled C cell) ---
$ws.Range("A4:I8").ClearContents()

# --- Rebuild hyperlinks: only keep C2 and C3 ---
# (Hyperlinks.Delete on any range clears every hyperlink in the sheet in this runtime)
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:omar@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Reem@gmail.com", "", "", "Reem@gmail.com")

# restore the original "Hyperlink" cell style (Add() re-applies a duplicate style entry)
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"

# set C3's displayed value (keeps the hyperlink's display text as the old email)
$ws.Range("C3").Value = "Gaber Mohamed"

# --- selection moves to C3 ---
$ws.Range("C3").Select()
